$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "59.763.67"
$ws.Range("E2").Value2 = "  +1.37%  "
$ws.Range("D3").Value2 = "2.299.00"
$ws.Range("E3").Value2 = "  -0.62%  "
$ws.Range("E4").Value2 = "  -0.02%  "
$ws.Range("D5").Value2 = "'539.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.00%  "
$ws.Range("D6").Value2 = "'129.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -2.37%  "
$ws.Range("E7").Value2 = "  -0.01%  "
$ws.Range("E8").Value2 = "  -2.24%  "
$ws.Range("D9").Value2 = "2.296.47"
$ws.Range("E9").Value2 = "  -0.66%  "
$ws.Range("D10").Value2 = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -0.55%  "
$ws.Range("D11").Value2 = "'5.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +0.39%  "
$ws.Range("E12").Value2 = "  -0.10%  "
$ws.Range("D13").Value2 = "'0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -1.05%  "
$ws.Range("D14").Value2 = "'23.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -2.77%  "
$ws.Range("B15").Value2 = "WrappedBTC"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value2 = "59.711.64"
$ws.Range("E15").Value2 = "  +1.49%  "
$ws.Range("B16").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value2 = "2.709.20"
$ws.Range("E16").Value2 = "  -0.59%  "
$ws.Range("D17").Value2 = "'0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -1.38%  "
$ws.Range("D18").Value2 = "2.311.55"
$ws.Range("E18").Value2 = "  -0.38%  "
$ws.Range("D19").Value2 = "'10.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -2.25%  "
$ws.Range("D20").Value2 = "'4.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -2.71%  "
$ws.Range("D21").Value2 = "'311.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -0.33%  "
$ws.Range("D22").Value2 = "'6.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -1.22%  "
$ws.Range("E23").Value2 = "  -0.30%  "
$ws.Range("D24").Value2 = "'5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.34%  "
$ws.Range("D25").Value2 = "'63.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +2.05%  "
$ws.Range("E26").Value2 = "  -2.90%  "
$ws.Range("E27").Value2 = "  -0.05%  "
$ws.Range("D28").Value2 = "'7.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -2.64%  "
$ws.Range("D29").Value2 = "'1.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +2.18%  "
$ws.Range("D30").Value2 = "'170.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.00%  "
$ws.Range("B31").Value2 = "SuiNetwork"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").Value2 = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +1.31%  "
$ws.Range("B32").Value2 = "PancakeSwap"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value2 = "'1.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -1.50%  "
$ws.Range("D33").Value2 = "0.0₃0721"
$ws.Range("E33").Value2 = "  -2.26%  "
$ws.Range("D34").Value2 = "'5.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -1.91%  "
$ws.Range("E35").Value2 = "  +1.49%  "
$ws.Range("D36").Value2 = "'0.377"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -2.61%  "
$ws.Range("E37").Value2 = "  +0.03%  "
$ws.Range("D38").Value2 = "'17.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -1.51%  "
$ws.Range("E39").Value2 = "  -0.02%  "
$ws.Range("E40").Value2 = "  -3.91%  "
$ws.Range("D41").Value2 = "'314.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +2.27%  "
$ws.Range("D42").Value2 = "'37.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -1.73%  "
$ws.Range("E43").Value2 = "  -0.71%  "
$ws.Range("D44").Value2 = "'135.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -4.17%  "
$ws.Range("D45").Value2 = "'3.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -1.28%  "
$ws.Range("D46").Value2 = "'0.0934"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -2.72%  "
$ws.Range("D47").Value2 = "'0.559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +0.33%  "
$ws.Range("D48").Value2 = "'18.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +1.43%  "
$ws.Range("D49").Value2 = "'0.0487"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -1.87%  "
$ws.Range("E50").Value2 = "  +16.15%  "
$ws.Range("D51").Value2 = "'0.0210"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -0.59%  "
